# dbo.xlsx: nueva tabla añadida
# Adds a new "lineas_fuentes_contactos_logs" table definition below the
# existing "lineas_fuentes_contactos" table (column F), rows 20-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table header (bold / shaded, same look as the other table headers) ---
# Copy the formatting of an existing header cell (F9) and paste it onto F20
# so the new header reuses the existing header style instead of creating a
# brand new one.
$ws.Range("F9").Copy()
$ws.Range("F20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F20").Value = "lineas_fuentes_contactos_logs"

# --- New table columns ---
$ws.Range("F21").Value = "id_linea_fuente_contacto_log"
$ws.Range("F22").Value = "id_linea_fuente"
$ws.Range("F23").Value = "accion"
$ws.Range("F24").Value = "resultado"

# Clear the clipboard marquee left over from the copy operation.
$excel.CutCopyMode = 0

# Match the saved cursor/selection position recorded for this edit.
[void]$ws.Range("H19").Select()
